# Auto-generated Excel COM-interop script
# Applies a market-price data refresh (scheduled runner) to the Leve profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 708583.0600000001
$ws.Range("J17").Value = 708583.0600000001
$ws.Range("L17").Value = 2125749.18
$ws.Range("N17").Value = -2126085.18
$ws.Range("H121").Value = 1074.975
$ws.Range("I121").Value = 880
$ws.Range("J121").Value = 1079.9744
$ws.Range("K121").Value = 2640
$ws.Range("L121").Value = 3239.9232
$ws.Range("M121").Value = -893
$ws.Range("N121").Value = -6733.9232
$ws.Range("H129").Value = 767.2258
$ws.Range("I129").Value = 501.26315
$ws.Range("J129").Value = 1188.3334
$ws.Range("K129").Value = 1503.78945
$ws.Range("L129").Value = 3565.0002
$ws.Range("M129").Value = 3496.21055
$ws.Range("N129").Value = -13565.0002
$ws.Range("H138").Value = 1930.05
$ws.Range("I138").Value = 1337.5625
$ws.Range("J138").Value = 2145.5
$ws.Range("K138").Value = 4012.6875
$ws.Range("L138").Value = 6436.5
$ws.Range("M138").Value = 1127.3125
$ws.Range("N138").Value = -16716.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15794262
$ws.Range("I32").Value = 1401393.5
$ws.Range("K32").Value = 1401393.5
$ws.Range("M32").Value = -1401106.5
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H97").Value = 538.4091
$ws.Range("I97").Value = 296.7647
$ws.Range("J97").Value = 1360
$ws.Range("K97").Value = 296.7647
$ws.Range("L97").Value = 1360
$ws.Range("M97").Value = 199.2353
$ws.Range("N97").Value = -2352

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1453.6072
$ws.Range("I86").Value = 1486.35
$ws.Range("J86").Value = 1371.75
$ws.Range("K86").Value = 1486.35
$ws.Range("L86").Value = 1371.75
$ws.Range("M86").Value = -363.3499999999999
$ws.Range("N86").Value = -3617.75
$ws.Range("H89").Value = 1453.6072
$ws.Range("I89").Value = 1486.35
$ws.Range("J89").Value = 1371.75
$ws.Range("K89").Value = 7431.75
$ws.Range("L89").Value = 6858.75
$ws.Range("M89").Value = -1815.75
$ws.Range("N89").Value = -18090.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 877.17
$ws.Range("I31").Value = 882.86957
$ws.Range("J31").Value = 872.3148
$ws.Range("K31").Value = 882.86957
$ws.Range("L31").Value = 872.3148
$ws.Range("M31").Value = -587.86957
$ws.Range("N31").Value = -1462.3148
$ws.Range("H34").Value = 877.17
$ws.Range("I34").Value = 882.86957
$ws.Range("J34").Value = 872.3148
$ws.Range("K34").Value = 882.86957
$ws.Range("L34").Value = 872.3148
$ws.Range("M34").Value = -680.86957
$ws.Range("N34").Value = -1276.3148
$ws.Range("H58").Value = 4469379
$ws.Range("I58").Value = 6804111.5
$ws.Range("J58").Value = 12163.546
$ws.Range("K58").Value = 6804111.5
$ws.Range("L58").Value = 12163.546
$ws.Range("M58").Value = -6803908.5
$ws.Range("N58").Value = -12569.546
$ws.Range("H132").Value = 7093980
$ws.Range("I132").Value = 9010204
$ws.Range("J132").Value = 3951
$ws.Range("K132").Value = 27030612
$ws.Range("L132").Value = 11853
$ws.Range("M132").Value = -27028082
$ws.Range("N132").Value = -16913
$ws.Range("H134").Value = 20834610
$ws.Range("I134").Value = 31251194
$ws.Range("J134").Value = 2718809.2
$ws.Range("K134").Value = 93753582
$ws.Range("L134").Value = 8156427.600000001
$ws.Range("M134").Value = -93751047
$ws.Range("N134").Value = -8161497.600000001
$ws.Range("H136").Value = 4469379
$ws.Range("I136").Value = 6804111.5
$ws.Range("J136").Value = 12163.546
$ws.Range("K136").Value = 20412334.5
$ws.Range("L136").Value = 36490.638
$ws.Range("M136").Value = -20409784.5
$ws.Range("N136").Value = -41590.638

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 1857
$ws.Range("I62").Value = 2000
$ws.Range("J62").Value = 1499.5
$ws.Range("K62").Value = 6000
$ws.Range("L62").Value = 4498.5
$ws.Range("M62").Value = -5314
$ws.Range("N62").Value = -5870.5
$ws.Range("H65").Value = 1857
$ws.Range("I65").Value = 2000
$ws.Range("J65").Value = 1499.5
$ws.Range("K65").Value = 18000
$ws.Range("L65").Value = 13495.5
$ws.Range("M65").Value = -14568
$ws.Range("N65").Value = -20359.5
$ws.Range("H68").Value = 2720.5781
$ws.Range("I68").Value = 610
$ws.Range("J68").Value = 4362.1387
$ws.Range("K68").Value = 1830
$ws.Range("L68").Value = 13086.4161
$ws.Range("M68").Value = -1019
$ws.Range("N68").Value = -14708.4161
$ws.Range("H71").Value = 2720.5781
$ws.Range("I71").Value = 610
$ws.Range("J71").Value = 4362.1387
$ws.Range("K71").Value = 5490
$ws.Range("L71").Value = 39259.24830000001
$ws.Range("M71").Value = -1434
$ws.Range("N71").Value = -47371.24830000001
$ws.Range("H80").Value = 4970.4116
$ws.Range("I80").Value = 1500
$ws.Range("J80").Value = 5187.3125
$ws.Range("K80").Value = 4500
$ws.Range("L80").Value = 15561.9375
$ws.Range("M80").Value = -3564
$ws.Range("N80").Value = -17433.9375
$ws.Range("H82").Value = 4250
$ws.Range("I82").Value = 2000
$ws.Range("J82").Value = 4571.4287
$ws.Range("K82").Value = 6000
$ws.Range("L82").Value = 13714.2861
$ws.Range("M82").Value = -5594
$ws.Range("N82").Value = -14526.2861
$ws.Range("H83").Value = 4970.4116
$ws.Range("I83").Value = 1500
$ws.Range("J83").Value = 5187.3125
$ws.Range("K83").Value = 13500
$ws.Range("L83").Value = 46685.8125
$ws.Range("M83").Value = -8820
$ws.Range("N83").Value = -56045.8125
$ws.Range("H85").Value = 4250
$ws.Range("I85").Value = 2000
$ws.Range("J85").Value = 4571.4287
$ws.Range("K85").Value = 6000
$ws.Range("L85").Value = 13714.2861
$ws.Range("M85").Value = -4596
$ws.Range("N85").Value = -16522.2861
$ws.Range("H97").Value = 733.4
$ws.Range("I97").Value = 295
$ws.Range("J97").Value = 1025.6666
$ws.Range("K97").Value = 885
$ws.Range("L97").Value = 3076.9998
$ws.Range("M97").Value = -389
$ws.Range("N97").Value = -4068.9998
$ws.Range("H98").Value = 1848
$ws.Range("I98").Value = 2595.7144
$ws.Range("J98").Value = 1324.6
$ws.Range("K98").Value = 7787.1432
$ws.Range("L98").Value = 3973.8
$ws.Range("M98").Value = -6289.1432
$ws.Range("N98").Value = -6969.799999999999
$ws.Range("H107").Value = 509.6771
$ws.Range("I107").Value = 232.77647
$ws.Range("J107").Value = 2649.3635
$ws.Range("K107").Value = 698.3294099999999
$ws.Range("L107").Value = 7948.0905
$ws.Range("M107").Value = 1221.67059
$ws.Range("N107").Value = -11788.0905
$ws.Range("H134").Value = 3237.1428
$ws.Range("I134").Value = 2332.3076
$ws.Range("J134").Value = 15000
$ws.Range("K134").Value = 6996.9228
$ws.Range("L134").Value = 45000
$ws.Range("M134").Value = -1926.9228
$ws.Range("N134").Value = -55140

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 32637.2
$ws.Range("I70").Value = 51562.668
$ws.Range("J70").Value = 4249
$ws.Range("K70").Value = 51562.668
$ws.Range("L70").Value = 4249
$ws.Range("M70").Value = -51292.668
$ws.Range("N70").Value = -4789
$ws.Range("H73").Value = 32637.2
$ws.Range("I73").Value = 51562.668
$ws.Range("J73").Value = 4249
$ws.Range("K73").Value = 51562.668
$ws.Range("L73").Value = 4249
$ws.Range("M73").Value = -50626.668
$ws.Range("N73").Value = -6121
$ws.Range("H113").Value = 1579.625
$ws.Range("I113").Value = 642.2
$ws.Range("J113").Value = 3142
$ws.Range("K113").Value = 642.2
$ws.Range("L113").Value = 3142
$ws.Range("M113").Value = 1527.8
$ws.Range("N113").Value = -7482

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 3727.6128
$ws.Range("I55").Value = 556.26086
$ws.Range("J55").Value = 12845.25
$ws.Range("K55").Value = 556.26086
$ws.Range("L55").Value = 12845.25
$ws.Range("M55").Value = -383.26086
$ws.Range("N55").Value = -13191.25
$ws.Range("H132").Value = 22227822
$ws.Range("I132").Value = 55561556
$ws.Range("J132").Value = 5332.6665
$ws.Range("K132").Value = 166684668
$ws.Range("L132").Value = 15997.9995
$ws.Range("M132").Value = -166682138
$ws.Range("N132").Value = -21057.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 25000
$ws.Range("I42").Value = 20000
$ws.Range("K42").Value = 20000
$ws.Range("M42").Value = -19622
$ws.Range("H80").Value = 18800
$ws.Range("I80").Value = 18800
$ws.Range("K80").Value = 18800
$ws.Range("M80").Value = -17802
$ws.Range("H83").Value = 18800
$ws.Range("I83").Value = 18800
$ws.Range("K83").Value = 56400
$ws.Range("M83").Value = -51408
$ws.Range("H132").Value = 81990130
$ws.Range("I132").Value = 137502820
$ws.Range("J132").Value = 1244408.1
$ws.Range("K132").Value = 412508460
$ws.Range("L132").Value = 3733224.3
$ws.Range("M132").Value = -412505930
$ws.Range("N132").Value = -3738284.3
$ws.Range("H136").Value = 32847852
$ws.Range("I136").Value = 31067498
$ws.Range("J136").Value = 38462812
$ws.Range("K136").Value = 93202494
$ws.Range("L136").Value = 115388436
$ws.Range("M136").Value = -93199944
$ws.Range("N136").Value = -115393536

Write-Host "Applied 249 cell updates across 8 sheets."
